$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (CO1, 30%)
$ws.Range("B8").Value = "Mothur"
$ws.Range("D8").Value = 0.47
$ws.Range("E8").Value = 0.74
$ws.Range("F8").Value = 0.5600000000000001
$ws.Range("G8").Value = 0.64
$ws.Range("H8").Value = 0.7

# Row 9 (CO1, 50%)
$ws.Range("B9").Value = "Mothur"
$ws.Range("D9").Value = 0.34
$ws.Range("E9").Value = 0.72
$ws.Range("F9").Value = 0.39
$ws.Range("G9").Value = 0.51
$ws.Range("H9").Value = 0.62

# Row 10 (CO1, 70%)
$ws.Range("B10").Value = "Mothur"
$ws.Range("D10").Value = 0.18
$ws.Range("E10").Value = 0.6899999999999999
$ws.Range("F10").Value = 0.2
$ws.Range("G10").Value = 0.31
$ws.Range("H10").Value = 0.46
